# #86 scene move cell implememnt
# Adds a new "TilePath" column (H) to the Scene table: header + type row +
# group row + a "default" value for every data row, then moves the
# selection to the last-edited cell (H9), matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1) Write the new header cell first so the table resize below can adopt it.
$ws.Range("H1").Value = "TilePath"

# 2) Grow the table (and its AutoFilter) from A1:G20 to A1:H20.
$tbl.Resize($ws.Range("A1:H20"))

# 3) Re-affirm the header text through the ListColumn's own range so the
#    table definition (xl/tables/table1.xml) picks up the real column name
#    instead of the auto-generated "Column8".
$lc = $tbl.ListColumns.Item($tbl.ListColumns.Count)
$lc.Range.Cells.Item(1).Value = "TilePath"

# 4) Copy the header-row / type-row formatting across so the new column
#    matches the look of the existing columns (bold header style, etc.).
$ws.Range("H2").Value = "string"
$ws.Range("H3").Value = "配置"

$ws.Range("A2").Copy()
$ws.Range("H2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A3").Copy()
$ws.Range("H3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("H2").Value = "string"
$ws.Range("H3").Value = "配置"

# 5) Fill every data row (4-20) with the "default" tile path value.
for ($r = 4; $r -le 20; $r++) {
    $ws.Range("H$r").Value = "default"
}

$excel.CutCopyMode = 0

# 6) Match the author's final selection.
[void]$ws.Range("H9").Select()
